$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '97.557.14'
$ws.Range('E2').Value = '  +1.72%  '
$ws.Range('D3').Value = '3.713.28'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D4').Value = '''2.36'
$ws.Range('E4').Value = '  +27.19%  '
$ws.Range('D5').Value = '''1.00'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = '''234.73'
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('D7').Value = '''661.19'
$ws.Range('E7').Value = '  +2.18%  '
$ws.Range('D8').Value = '''0.452'
$ws.Range('E8').Value = '  +8.22%  '
$ws.Range('D9').Value = '''1.15'
$ws.Range('E9').Value = '  +10.42%  '
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').Value = '3.709.83'
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').Value = '''45.44'
$ws.Range('E12').Value = '  +3.55%  '
$ws.Range('D13').Value = '''0.0000311'
$ws.Range('E13').Value = '  +11.18%  '
$ws.Range('D14').Value = '''0.211'
$ws.Range('E14').Value = '  +4.06%  '
$ws.Range('D15').Value = '''6.78'
$ws.Range('E15').Value = '  +1.63%  '
$ws.Range('D16').Value = '4.424.02'
$ws.Range('E16').Value = '  +1.55%  '
$ws.Range('D17').Value = '97.636.48'
$ws.Range('E17').Value = '  +1.91%  '
$ws.Range('D18').Value = '''8.97'
$ws.Range('E18').Value = '  +15.28%  '
$ws.Range('D19').Value = '3.716.74'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '''19.22'
$ws.Range('E20').Value = '  +3.84%  '
$ws.Range('D21').Value = '''13.13'
$ws.Range('E21').Value = '  +3.61%  '
$ws.Range('D22').Value = '''0.549'
$ws.Range('E22').Value = '  +10.73%  '
$ws.Range('D23').Value = '''546.38'
$ws.Range('E23').Value = '  +6.72%  '
$ws.Range('D24').Value = '''3.40'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('B25').Value = 'PEPE'
$ws.Range('C25').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D25').Value = '''0.0000221'
$ws.Range('E25').Value = '  +9.43%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '''123.10'
$ws.Range('E26').Value = '  +22.70%  '
$ws.Range('E27').Value = '  +42.84%  '
$ws.Range('D28').Value = '''6.87'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('B29').Value = 'WrappedeETH'
$ws.Range('C29').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D29').Value = '3.912.51'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '''13.29'
$ws.Range('E30').Value = '  +2.80%  '
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').Value = '''13.19'
$ws.Range('E31').Value = '  +10.13%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''3.09'
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('B33').Value = 'Dai'
$ws.Range('C33').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D33').Value = '''0.999'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').Value = '''0.187'
$ws.Range('E34').Value = '  +3.27%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '''33.65'
$ws.Range('E35').Value = '  +5.95%  '
$ws.Range('B36').Value = 'Fetch.AI'
$ws.Range('C36').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D36').Value = '''1.84'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('B37').Value = 'PolygonEcosystemToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D37').Value = '''0.622'
$ws.Range('E37').Value = '  +7.16%  '
$ws.Range('B38').Value = 'Binance-PegBSC-USD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D38').Value = '''0.997'
$ws.Range('E38').Value = '  -0.12%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '''631.37'
$ws.Range('E39').Value = '  -3.14%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').Value = '''8.57'
$ws.Range('E41').Value = '  -1.10%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''7.18'
$ws.Range('E42').Value = '  +7.43%  '
$ws.Range('B43').Value = 'Kaspa'
$ws.Range('C43').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D43').Value = '''0.168'
$ws.Range('E43').Value = '  +6.79%  '
$ws.Range('D44').Value = '''2.05'
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0502'
$ws.Range('E45').Value = '  +13.49%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''40.34'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '''0.485'
$ws.Range('E47').Value = '  +12.54%  '
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').Value = '''0.970'
$ws.Range('E48').Value = '  +2.57%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '''2.39'
$ws.Range('E49').Value = '  +6.86%  '
$ws.Range('B50').Value = 'Cosmos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D50').Value = '''9.09'
$ws.Range('E50').Value = '  +9.11%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '''23.61'
$ws.Range('E51').Value = '  +0.36%  '
